$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A2 value changes from "Apportionment" to a new (typo) string "Approtionment"
$ws.Range("A2").Value = "Approtionment"

# Duplicate existing rows 2 and 3 into new rows 6 and 7 (keeps styles/number formats intact)
$ws.Range("A2:X2").Copy($ws.Range("A6:X6"))
$ws.Range("A3:X3").Copy($ws.Range("A7:X7"))

# Row 6 specific differences vs row 2
$ws.Range("B6").Value = "Demo 7"
$ws.Range("W6").Value = "Decrease"

# Row 7 specific differences vs row 3
$ws.Range("B7").Value = "Demo 7"

# Update the view: select W7 (becomes the active cell / selection)
$ws.Activate()
$ws.Range("W7").Select()
